# Auto-generated edit script: updates cryptos list values (price/volume) and
# swaps the Aave / RenderToken rows (45 <-> 46) per the Oct 12 2023 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''26.857.88'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '''  -1.12%  '
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = '''1.563.71'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '''  +0.05%  '
$ws.Range("E3").ClearFormats()
$ws.Range("E4").Value = '''  -0.09%  '
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = '''206.07'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '''  -0.30%  '
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = '''0.488'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '''  -1.17%  '
$ws.Range("E6").ClearFormats()
$ws.Range("E7").Value = '''  -0.05%  '
$ws.Range("E7").ClearFormats()
$ws.Range("E8").Value = '''  -2.16%  '
$ws.Range("E8").ClearFormats()
$ws.Range("E9").Value = '''  -0.45%  '
$ws.Range("E9").ClearFormats()
$ws.Range("E10").Value = '''  -1.19%  '
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = '''0.0863'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '''  +0.22%  '
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = '''1.785.89'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '''  +0.06%  '
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = '''1.570.89'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '''  +0.38%  '
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = '''3.72'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '''  -1.28%  '
$ws.Range("E14").ClearFormats()
$ws.Range("E15").Value = '''  -0.11%  '
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = '''26.867.61'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '''  -1.05%  '
$ws.Range("E16").ClearFormats()
$ws.Range("E17").Value = '''  -2.68%  '
$ws.Range("E17").ClearFormats()
$ws.Range("E18").Value = '''  +0.83%  '
$ws.Range("E18").ClearFormats()
$ws.Range("E19").Value = '''  +1.85%  '
$ws.Range("E19").ClearFormats()
$ws.Range("E20").Value = '''  -1.17%  '
$ws.Range("E20").ClearFormats()
$ws.Range("E22").Value = '''  +0.24%  '
$ws.Range("E22").ClearFormats()
$ws.Range("E23").Value = '''  -2.13%  '
$ws.Range("E23").ClearFormats()
$ws.Range("E24").Value = '''  +1.41%  '
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = '''153.41'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '''  +0.91%  '
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = '''6.73'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '''  +2.45%  '
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = '''14.90'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '''  +0.12%  '
$ws.Range("E27").ClearFormats()
$ws.Range("E28").Value = '''  -0.09%  '
$ws.Range("E28").ClearFormats()
$ws.Range("E29").Value = '''  -0.89%  '
$ws.Range("E29").ClearFormats()
$ws.Range("D30").Value = '''0.0466'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '''  +0.40%  '
$ws.Range("E30").ClearFormats()
$ws.Range("E31").Value = '''  -3.41%  '
$ws.Range("E31").ClearFormats()
$ws.Range("D32").Value = '''3.16'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '''  -0.07%  '
$ws.Range("E32").ClearFormats()
$ws.Range("D33").Value = '''1.403.91'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '''  +1.51%  '
$ws.Range("E33").ClearFormats()
$ws.Range("E34").Value = '''  -0.75%  '
$ws.Range("E34").ClearFormats()
$ws.Range("E35").Value = '''  -1.47%  '
$ws.Range("E35").ClearFormats()
$ws.Range("E36").Value = '''  -0.46%  '
$ws.Range("E36").ClearFormats()
$ws.Range("D37").Value = '''0.919'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '''  -2.20%  '
$ws.Range("E37").ClearFormats()
$ws.Range("E38").Value = '''  -0.46%  '
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = '''0.529'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '''  +2.30%  '
$ws.Range("E39").ClearFormats()
$ws.Range("E40").Value = '''  -0.53%  '
$ws.Range("E40").ClearFormats()
$ws.Range("E41").Value = '''  -0.06%  '
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = '''0.997'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '''  +0.56%  '
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = '''5.44'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '''  +4.12%  '
$ws.Range("E43").ClearFormats()
$ws.Range("E44").Value = '''  +0.38%  '
$ws.Range("E44").ClearFormats()
$ws.Range("B45").Value = '''RenderToken'
$ws.Range("B45").ClearFormats()
$ws.Range("C45").Value = '''https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("C45").ClearFormats()
$ws.Range("D45").Value = '''1.77'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '''  -0.92%  '
$ws.Range("E45").ClearFormats()
$ws.Range("B46").Value = '''Aave'
$ws.Range("B46").ClearFormats()
$ws.Range("C46").Value = '''https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("C46").ClearFormats()
$ws.Range("D46").Value = '''63.42'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '''  -0.02%  '
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = '''1.699.34'
$ws.Range("D47").ClearFormats()
$ws.Range("D48").Value = '''86.63'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '''  +1.19%  '
$ws.Range("E48").ClearFormats()
$ws.Range("E49").Value = '''  +2.95%  '
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = '''0.0₇0974'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '''  -2.10%  '
$ws.Range("E50").ClearFormats()
$ws.Range("E51").Value = '''  +0.75%  '
$ws.Range("E51").ClearFormats()
